# Insert a new row before row 47 (shifts old rows 47-98 down to 48-99)
# and fill the new row 47 with the data for the new weekly entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(47).Insert()

$ws.Cells.Item(47, 1).Value = 8
$ws.Cells.Item(47, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(47, 3).Value = "Coquimbo"
$ws.Cells.Item(47, 4).Value = 44512
$ws.Cells.Item(47, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(47, 5).Value = 4
$ws.Cells.Item(47, 6).Value = 100112001
$ws.Cells.Item(47, 7).Value = "Berenjena"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 600
$ws.Cells.Item(47, 11).Value = 8000
$ws.Cells.Item(47, 12).Value = 9000
$ws.Cells.Item(47, 13).Value = 8500
$ws.Cells.Item(47, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(47, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(47, 16).Value = 142
$ws.Cells.Item(47, 17).Value = 60
$ws.Cells.Item(47, 18).Value = "Hortaliza"
